# Time_spent.xlsx — "Add missing testing output, fix diagrams and time chart"
#
# Updates the monthly time-tracking figures for the "Analýza" (row 2),
# "Návrh" (row 3) and "Implementace" (row 4) series, nudges the embedded
# line-chart a touch to the right, and leaves the selection on F4 (where
# the author's cursor ended up after editing that cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Analýza ------------------------------------------------------
$ws.Range("C2").Value2 = 0.61597222222222225
$ws.Range("E2").Value2 = 0.25
$ws.Range("F2").Value2 = 0.16666666666666666
$ws.Range("G2").Value2 = 0.083333333333333329

# --- Row 3: Návrh ---------------------------------------------------------
$ws.Range("C3").Value2 = 0.25
$ws.Range("E3").Value2 = 0.22013888888888888
$ws.Range("F3").Value2 = 0.30208333333333331

# --- Row 4: Implementace --------------------------------------------------
$ws.Range("F4").Value2 = 0.22847222222222222
$ws.Range("G4").Value2 = 2.3666666666666667

# --- Nudge the time-spent chart 2.1pt (26670 EMU) to the right -----------
$co = $ws.ChartObjects(1)
$co.Left = $co.Left + 2.1

# --- Final selection / view state -----------------------------------------
$ws.Range("F4").Select() | Out-Null
